# translater traverses folder, writes generated tokens back to files
#
# Sheet1 keeps its "hun" column (A) as-is, but the "eng" column (B) is
# dropped. Three new sheets are added: a second source-language sheet
# ("Sheet2") and a translated companion for each source sheet
# ("Translated_Sheet1" / "Translated_Sheet2") with a bold, bordered,
# centered header row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: drop the now-unused "eng" header/column, keep column A intact.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sheet1"
$ws1.Range("B1").ClearContents()

# ---------------------------------------------------------------------
# Sheet2: a second small source sheet (hun words/phrases).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$ws2.Cells.Item(1, 1).Value = "hun"
$ws2.Cells.Item(2, 1).Value = "körte"
$ws2.Cells.Item(3, 1).Value = "szék"
$ws2.Cells.Item(4, 1).Value = "Ez egy egész mondat."

$ws2.Columns.Item(1).ColumnWidth = 17.1
$ws2.Range("B6").Select()

# ---------------------------------------------------------------------
# Translated_Sheet1: hun/eng pairs that mirror Sheet1.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "Translated_Sheet1"

$ws3.Cells.Item(1, 1).Value = "hun"
$ws3.Cells.Item(1, 2).Value = "translated_hun"
$ws3.Cells.Item(2, 1).Value = "alma"
$ws3.Cells.Item(2, 2).Value = "['Apples']"
$ws3.Cells.Item(3, 1).Value = "nem gondolom, hogy ez nehezére esne"
$ws3.Cells.Item(3, 2).Value = "[""I don't think that's going to be difficult.""]"
$ws3.Cells.Item(4, 1).Value = "majd meglátjuk"
$ws3.Cells.Item(4, 2).Value = "[""We'll see""]"
$ws3.Cells.Item(5, 1).Value = "panzerkraftwagen"
$ws3.Cells.Item(5, 2).Value = "['Other vehicles']"

$headerRange3 = $ws3.Range("A1:B1")
$headerRange3.Font.Bold = $true
$headerRange3.Borders.LineStyle = 1
$headerRange3.HorizontalAlignment = -4108
$headerRange3.VerticalAlignment = -4160

$ws3.Columns.Item(2).ColumnWidth = 32.92

# ---------------------------------------------------------------------
# Translated_Sheet2: hun/eng pairs that mirror Sheet2.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws3)
$ws4.Name = "Translated_Sheet2"

$ws4.Cells.Item(1, 1).Value = "hun"
$ws4.Cells.Item(1, 2).Value = "translated_hun"
$ws4.Cells.Item(2, 1).Value = "körte"
$ws4.Cells.Item(2, 2).Value = "['the roasting']"
$ws4.Cells.Item(3, 1).Value = "szék"
$ws4.Cells.Item(3, 2).Value = "['Other, of a kind used for the manufacture of goods']"
$ws4.Cells.Item(4, 1).Value = "Ez egy egész mondat."
$ws4.Cells.Item(4, 2).Value = "[""That's a whole sentence.""]"

$headerRange4 = $ws4.Range("A1:B1")
$headerRange4.Font.Bold = $true
$headerRange4.Borders.LineStyle = 1
$headerRange4.HorizontalAlignment = -4108
$headerRange4.VerticalAlignment = -4160

$ws4.Columns.Item(1).ColumnWidth = 24.75
$ws4.Columns.Item(2).ColumnWidth = 42.26

# Leave Sheet1 as the active/selected sheet, matching the source file.
$ws1.Select()
